$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 data
$ws.Range("D2").Value = "Racing"
$ws.Range("B2").Value = "tafur.fredy@gmail.com"
$ws.Range("A2").Value = "Fredy Tafur Garay"
$ws.Range("C2").Value = 41863284
# F2 stays "socio"

# Clear row 3 data (A3, C3, D3, F3), keep B3 cell present but empty
$ws.Range("A3").ClearContents()
$ws.Range("C3").ClearContents()
$ws.Range("D3").ClearContents()
$ws.Range("F3").ClearContents()
$ws.Range("B3").ClearContents()

# Remove hyperlinks from the sheet (keeps cell values/styles)
$ws.Hyperlinks.Delete()

# Update selection
$ws.Range("E6").Select()
